$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Truth table edits: row 4 (boz) updates
$ws.Range("B4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("N4").Value = 0

# Update the active selection to match the saved view state
$ws.Range("N4").Select()
